# Generate Report for Handback
#
# Marks the 6a0f19e4... and 83b83e2c... files as handed back (in sync with
# en-US) on every sheet, records the handback timestamp in the "Latest
# Handback DateTime" column, and fills in the "Latest Target File" /
# "Latest Handback File" columns (with working hyperlinks) on both the
# zh-cn and de-de detail sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: just the rolled-up status column per language
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $statusText
$ov.Range("C2").Value = $statusText
$ov.Range("B3").Value = $statusText
$ov.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = $statusText
$zh.Range("B3").Value = $statusText

# Latest Handback DateTime
$zh.Range("G2").Value = "2016-03-09 01:05:24"
$zh.Range("G3").Value = "2016-03-09 01:05:24"

# Latest Target File (E) / Latest Handback File (F) for row 2
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/3613c5c9ea17bc3d31e0c59c12dce6411dca7545/e2e/6a0f19e4-5607-4441-ace2-f1b6fd638a06.md", [Type]::Missing, [Type]::Missing, "6a0f19e4-5607-4441-ace2-f1b6fd638a06.md")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/29551833dd2aa03d3f9a964402b6bfd96aa91ef8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6a0f19e4-5607-4441-ace2-f1b6fd638a06.c1fb9592cb7e8835f8b89d21457f65f66582be87.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "6a0f19e4-5607-4441-ace2-f1b6fd638a06.c1fb9592cb7e8835f8b89d21457f65f66582be87.zh-cn.xlf")

# Latest Target File (E) / Latest Handback File (F) for row 3
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/3613c5c9ea17bc3d31e0c59c12dce6411dca7545/e2e/83b83e2c-b8da-45f8-b9d7-a277efb4ec8e.md", [Type]::Missing, [Type]::Missing, "83b83e2c-b8da-45f8-b9d7-a277efb4ec8e.md")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/29551833dd2aa03d3f9a964402b6bfd96aa91ef8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/83b83e2c-b8da-45f8-b9d7-a277efb4ec8e.983c52d35c0c8530653e2e8366471960e7b12c4e.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "83b83e2c-b8da-45f8-b9d7-a277efb4ec8e.983c52d35c0c8530653e2e8366471960e7b12c4e.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = $statusText
$de.Range("B3").Value = $statusText

# Latest Handback DateTime
$de.Range("G2").Value = "2016-03-09 01:05:55"
$de.Range("G3").Value = "2016-03-09 01:05:55"

# Latest Target File (E) / Latest Handback File (F) for row 2
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/3613c5c9ea17bc3d31e0c59c12dce6411dca7545/e2e/6a0f19e4-5607-4441-ace2-f1b6fd638a06.md", [Type]::Missing, [Type]::Missing, "6a0f19e4-5607-4441-ace2-f1b6fd638a06.md")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef63039aa844c045ca43e67b69f19e62419a3d25/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6a0f19e4-5607-4441-ace2-f1b6fd638a06.c1fb9592cb7e8835f8b89d21457f65f66582be87.de-de.xlf", [Type]::Missing, [Type]::Missing, "6a0f19e4-5607-4441-ace2-f1b6fd638a06.c1fb9592cb7e8835f8b89d21457f65f66582be87.de-de.xlf")

# Latest Target File (E) / Latest Handback File (F) for row 3
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/3613c5c9ea17bc3d31e0c59c12dce6411dca7545/e2e/83b83e2c-b8da-45f8-b9d7-a277efb4ec8e.md", [Type]::Missing, [Type]::Missing, "83b83e2c-b8da-45f8-b9d7-a277efb4ec8e.md")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef63039aa844c045ca43e67b69f19e62419a3d25/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/83b83e2c-b8da-45f8-b9d7-a277efb4ec8e.983c52d35c0c8530653e2e8366471960e7b12c4e.de-de.xlf", [Type]::Missing, [Type]::Missing, "83b83e2c-b8da-45f8-b9d7-a277efb4ec8e.983c52d35c0c8530653e2e8366471960e7b12c4e.de-de.xlf")
